$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Update column O (faturamento) values for rows 2-6
$ws.Range("O2").Value = 14329.8
$ws.Range("O3").Value = 2756.9
$ws.Range("O4").Value = 2421
$ws.Range("O5").Value = 1726.26
$ws.Range("O6").Value = 21233.96

# Update column AG (total) values for rows 2-6 to reflect new totals
$ws.Range("AG2").Value = 151838.88
$ws.Range("AG3").Value = 72874.3
$ws.Range("AG4").Value = 40327.15
$ws.Range("AG5").Value = 37758.62
$ws.Range("AG6").Value = 302798.95
